$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text representation
# (Excel auto-converts plain-numeric-looking strings to numbers on assignment,
# which would lose trailing zeros / exact formatting). Force text format first,
# then restore the default "Normal" style so no stray formatting is introduced.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '66.528.41'
$ws.Range('D3').Value = '2.504.59'
$ws.Range('E3').Value = '  -4.88%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '582.70'
$ws.Range('E5').Value = '  -2.09%  '
$ws.Range('D6').Value = '172.19'
$ws.Range('E6').Value = '  +2.85%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -1.65%  '
$ws.Range('D9').Value = '2.503.99'
$ws.Range('E9').Value = '  -4.88%  '
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('E12').Value = '  -3.29%  '
$ws.Range('E13').Value = '  -2.23%  '
$ws.Range('D14').Value = '26.61'
$ws.Range('E14').Value = '  -3.72%  '
$ws.Range('D15').Value = '2.969.95'
$ws.Range('D16').Value = '0.0000176'
$ws.Range('E16').Value = '  -2.78%  '
$ws.Range('D17').Value = '66.298.78'
$ws.Range('E17').Value = '  -1.70%  '
$ws.Range('D18').Value = '2.506.10'
$ws.Range('E18').Value = '  -4.65%  '
$ws.Range('E19').Value = '  -3.70%  '
$ws.Range('D20').Value = '11.24'
$ws.Range('E20').Value = '  -6.11%  '
$ws.Range('D21').Value = '347.86'
$ws.Range('E21').Value = '  -2.52%  '
$ws.Range('E22').Value = '  -2.71%  '
$ws.Range('E23').Value = '  -0.79%  '
$ws.Range('D24').Value = '1.98'
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').Value = '69.69'
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').Value = '10.03'
$ws.Range('E27').Value = '  -2.33%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('D29').Value = '2.640.32'
$ws.Range('E29').Value = '  -4.43%  '
$ws.Range('D30').Value = '0.0₃0979'
$ws.Range('E30').Value = '  -2.88%  '
$ws.Range('D31').Value = '529.52'
$ws.Range('E31').Value = '  -3.34%  '
$ws.Range('D32').Value = '8.08'
$ws.Range('E32').Value = '  +1.89%  '
$ws.Range('E33').Value = '  -2.41%  '
$ws.Range('E34').Value = '  -3.04%  '
$ws.Range('E35').Value = '  -3.99%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('E37').Value = '  -2.77%  '
$ws.Range('D38').Value = '156.58'
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('D39').Value = '18.61'
$ws.Range('E39').Value = '  -2.02%  '
$ws.Range('D40').Value = '18.36'
$ws.Range('E40').Value = '  +0.35%  '
$ws.Range('E41').Value = '  -3.29%  '
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('D43').Value = '5.10'
$ws.Range('D44').Value = '2.53'
$ws.Range('E44').Value = '  +4.20%  '
$ws.Range('D46').Value = '39.53'
$ws.Range('E46').Value = '  -1.37%  '
$ws.Range('D47').Value = '148.32'
$ws.Range('E47').Value = '  -2.99%  '
$ws.Range('E48').Value = '  -3.85%  '
$ws.Range('E49').Value = '  -3.49%  '
$ws.Range('D50').Value = '1.72'
$ws.Range('E50').Value = '  +1.35%  '
$ws.Range('E51').Value = '  -8.59%  '

$dataRange.Style = "Normal"
